$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 2) ------------------------------------------------
$ws.Range("A2").Value = "MCH193-1"
$ws.Range("C2").Value = "MRAP- NEWSLETTERS, MAGAZINES"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24A | GRAP COUNT NUMER: NONE"

# --- Formatting for the new row -------------------------------------------
# Build the base style (10pt Calibri, automatic/theme text colour) on A2,
# then fan it out to the other cells that share that exact style so they
# all resolve to the same cell-format record.
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.ThemeColor = 1

$ws.Range("A2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("F2").PasteSpecial(-4122)

# Re-apply the value overwritten by paste (PasteSpecial formats shouldn't
# touch values, but make sure F2 keeps its text either way).
$ws.Range("F2").Value = "1 Box"

$excel.CutCopyMode = $false

# --- Restore the frozen header pane / selection on the new row ------------
$ws.Range("A2:I2").Select()
$excel.ActiveWindow.FreezePanes = $true
